$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values in column E (STEP 8 + STEP 9 data correction)
$ws.Range("E3").Value = 216
$ws.Range("E4").Value = 72
$ws.Range("E5").Value = 144
$ws.Range("E6").Value = 72
$ws.Range("E7").Value = 72
$ws.Range("E8").Value = 216
$ws.Range("E9").Value = 216

$ws.Range("E11").Value = 25920
$ws.Range("E12").Value = 8640
$ws.Range("E13").Value = 17280
$ws.Range("E14").Value = 8640
$ws.Range("E15").Value = 8640
$ws.Range("E16").Value = 25920
$ws.Range("E17").Value = 25920

# Update the selected range on the sheet (entire column F selected)
$ws.Range("F1:F1048576").Select()
